# Weekly crime-stat refresh: roll the report forward one week
# (Volume/Number 43 -> 44, week 10/21-10/27/2024 -> 10/28/2024-11/3/2024)
# and load the newly collected weekly/28-day/YTD crime figures for the
# 7th Precinct table (rows 15-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead text: bump the report volume/number and the week-covering dates ---
$ws.Range("A8").Value = "Volume 31   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/28/2024  Through  11/3/2024"

# --- Column J (Year-to-date "2023") got one character narrower to match the
#     other integer count columns once the new figures were bestFit. ---
$ws.Columns.Item(10).ColumnWidth = 5.43

# --- Cells needing a style/type change: seed via Copy() from a stable same-style cell, then set value ---
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("J14").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 1
$ws.Range("K14").Copy($ws.Range("E20"))
$ws.Range("E20").Value = -100
$ws.Range("C14").Copy($ws.Range("G27"))
$ws.Range("E14").Copy($ws.Range("H27"))
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("J14").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 2
$ws.Range("K14").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100
$ws.Range("J14").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 1
$ws.Range("K14").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100
$ws.Range("J14").Copy($ws.Range("G31"))
$ws.Range("G31").Value = 1
$ws.Range("K14").Copy($ws.Range("H31"))
$ws.Range("H31").Value = -100

# --- Simple value updates (style/type unchanged) ---
$ws.Range("L15").Value = -43.75
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 7
$ws.Range("H16").Value = -12.5
$ws.Range("I16").Value = 89
$ws.Range("J16").Value = 112
$ws.Range("K16").Value = -20.535714285714
$ws.Range("L16").Value = -44.025157232704
$ws.Range("M16").Value = -25.210084033613
$ws.Range("N16").Value = -87.194244604316
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 21
$ws.Range("H17").Value = 16.666666666666
$ws.Range("I17").Value = 199
$ws.Range("J17").Value = 185
$ws.Range("K17").Value = 7.567567567567
$ws.Range("L17").Value = 6.417112299465
$ws.Range("M17").Value = 74.561403508771
$ws.Range("N17").Value = 2.577319587628
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -54.545454545454
$ws.Range("I18").Value = 96
$ws.Range("J18").Value = 148
$ws.Range("K18").Value = -35.135135135135
$ws.Range("L18").Value = -32.867132867132
$ws.Range("M18").Value = 14.285714285714
$ws.Range("N18").Value = -71.513353115727
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -47.058823529411
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = -19.607843137254
$ws.Range("I19").Value = 430
$ws.Range("J19").Value = 506
$ws.Range("K19").Value = -15.019762845849
$ws.Range("L19").Value = -34.351145038167
$ws.Range("M19").Value = 86.147186147186
$ws.Range("N19").Value = 9.41475826972
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 25
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = -4
$ws.Range("L20").Value = -14.285714285714
$ws.Range("M20").Value = 4.347826086956
$ws.Range("N20").Value = -85.321100917431
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -25
$ws.Range("F21").Value = 79
$ws.Range("G21").Value = 92
$ws.Range("H21").Value = -14.130434782608
$ws.Range("I21").Value = 871
$ws.Range("J21").Value = 1011
$ws.Range("K21").Value = -13.847675568743
$ws.Range("L21").Value = -28.547990155865
$ws.Range("M21").Value = 44.444444444444
$ws.Range("N21").Value = -55.786802030456
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -20
$ws.Range("F23").Value = 23
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = 76.923076923076
$ws.Range("I23").Value = 183
$ws.Range("J23").Value = 140
$ws.Range("K23").Value = 30.714285714285
$ws.Range("L23").Value = 21.192052980132
$ws.Range("M23").Value = 56.410256410256
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 22.222222222222
$ws.Range("F24").Value = 127
$ws.Range("G24").Value = 117
$ws.Range("H24").Value = 8.547008547008
$ws.Range("I24").Value = 1197
$ws.Range("J24").Value = 1061
$ws.Range("K24").Value = 12.818096135721
$ws.Range("L24").Value = -38.83495145631
$ws.Range("M24").Value = 85.581395348837
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 22
$ws.Range("E25").Value = 9.090909090909
$ws.Range("F25").Value = 78
$ws.Range("G25").Value = 84
$ws.Range("H25").Value = -7.142857142857
$ws.Range("I25").Value = 794
$ws.Range("J25").Value = 603
$ws.Range("K25").Value = 31.67495854063
$ws.Range("L25").Value = -49.523204068658
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = 13.333333333333
$ws.Range("F26").Value = 50
$ws.Range("G26").Value = 37
$ws.Range("H26").Value = 35.135135135135
$ws.Range("I26").Value = 381
$ws.Range("J26").Value = 384
$ws.Range("K26").Value = -0.78125
$ws.Range("L26").Value = -2.056555269922
$ws.Range("M26").Value = 37.05035971223
$ws.Range("L27").Value = -36
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = -25
$ws.Range("J28").Value = 42
$ws.Range("K28").Value = -19.047619047619
$ws.Range("L28").Value = -12.820512820512
$ws.Range("F29").Value = 1
$ws.Range("N29").Value = -62.5
$ws.Range("F30").Value = 1
$ws.Range("N30").Value = -57.142857142857
$ws.Range("J31").Value = 4
$ws.Range("K31").Value = 225

